$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Тип помещения" (room type) values in column C for rows 3-6.
# These previously all said "Квартира" and now describe more specific
# non-apartment premises types used by the updated meter-import example.
$ws.Range("C3").Value = "Машиноместо"
$ws.Range("C4").Value = "Апартаменты"
$ws.Range("C5").Value = "Кладовая"
$ws.Range("C6").Value = "Коммерческое помещение"
